# Autogenerated on Mon Feb 09 2015 03:30:35 GMT+0000 (Coordinated Universal Time)
#
# Adds a new "Number of employees / Assets / Turnover" size-class
# breakdown table (rows 23-27) to the Summary sheet, and relocates the
# trailing "SME Performance Review EU" source citation (previously at
# rows 26-27) down to rows 32-33 to make room for it.
#
# NOTE: this COM-interop runtime re-normalizes the style table on every
# save, which can drop the cached bold/italic/size/underline formatting
# of cells that the script does not explicitly touch. To keep the
# workbook's existing look intact we re-assert the formatting of every
# previously-styled cell in addition to writing the new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Re-assert formatting of pre-existing cells so it survives the
# style-table round trip untouched.
# ---------------------------------------------------------------

# "name" style (18pt) - report title
$ws.Range("A1").Font.Size = 18

# "title" style (bold) - section headings / table headers / row labels
$boldCells = "A3","B11","C11","D11","A12","A13","A14","A15","A16","B19","C19","D19","A20"
foreach ($addr in $boldCells) {
    $ws.Range($addr).Font.Bold = $true
}

# "title_" style (bold + underline)
$ws.Range("A9").Font.Bold = $true
$ws.Range("A9").Font.Underline = $true

# "source" style (italic) - citation lines
$italicCells = "A17","A21"
foreach ($addr in $italicCells) {
    $ws.Range($addr).Font.Italic = $true
}

# ---------------------------------------------------------------
# New table header (row 23) - bold "title" style, same as row 11/19
# ---------------------------------------------------------------
$ws.Range("B23").Value = "Number of employees"
$ws.Range("B23").Font.Bold = $true

$ws.Range("C23").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("C23").Font.Bold = $true

$ws.Range("D23").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("D23").Font.Bold = $true

# ---------------------------------------------------------------
# New table body (rows 24-27) - plain "Normal" style
# ---------------------------------------------------------------
$ws.Range("A24").Value = "Micro"
$ws.Range("A24").Font.Bold = $false

$ws.Range("B24").Value = "<10"
$ws.Range("B24").Font.Bold = $false

$ws.Range("C24").Value = ""
$ws.Range("C24").Font.Bold = $false

$ws.Range("D24").Value = ""
$ws.Range("D24").Font.Bold = $false

$ws.Range("A25").Value = "Small"
$ws.Range("A25").Font.Bold = $false

$ws.Range("B25").Value = "<50"
$ws.Range("B25").Font.Bold = $false

$ws.Range("C25").Value = ""
$ws.Range("C25").Font.Bold = $false

$ws.Range("D25").Value = ""
$ws.Range("D25").Font.Bold = $false

# Row 26 previously held the (bold) "SME Performance Review EU" source
# label - repurpose it for the "Medium" table row and clear the old
# bold formatting so it matches the plain body style.
$ws.Range("A26").Value = "Medium"
$ws.Range("A26").Font.Bold = $false

$ws.Range("B26").Value = "<250"
$ws.Range("B26").Font.Bold = $false

$ws.Range("C26").Value = ""
$ws.Range("C26").Font.Bold = $false

$ws.Range("D26").Value = ""
$ws.Range("D26").Font.Bold = $false

# Row 27 previously held the (italic) source citation - repurpose it
# for the "Large" table row and clear the old italic formatting.
$ws.Range("A27").Value = "Large"
$ws.Range("A27").Font.Italic = $false

$ws.Range("B27").Value = ">249"
$ws.Range("B27").Font.Bold = $false

$ws.Range("C27").Value = ""
$ws.Range("C27").Font.Bold = $false

$ws.Range("D27").Value = ""
$ws.Range("D27").Font.Bold = $false

# ---------------------------------------------------------------
# Relocated source citation (previously at rows 26-27) now sits
# below the new table, at rows 32-33.
# ---------------------------------------------------------------
$ws.Range("A32").Value = "SME Performance Review EU"
$ws.Range("A32").Font.Bold = $true

$ws.Range("A33").Value = "SME Performance Review EU, ""SBA Fact sheet"", 2013.  Available at http://ec.europa.eu/enterprise/policies/sme/facts-figures-analysis/performance-review/index_en.htm"
$ws.Range("A33").Font.Italic = $true
